# Update tab names in all BOMs, fix bi-color LED naming.

$wb = $excel.ActiveWorkbook

# Rename the worksheet tab from "WF" to the generic "BOM" label.
$ws = $wb.Worksheets.Item("WF")
$ws.Name = "BOM"

# The Ref cells for the bi-color LED trim pots / headers (D22:D25) carried an
# explicit (redundant) cell style; clear it back to the workbook's default
# "Normal" formatting so they match the rest of the Ref column.
$ws.Range("D22:D25").Style = "Normal"
